$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 12
$ws.Range("D7").Value = 5

$ws.Range("B8").Value = 18
$ws.Range("D8").Value = 12

$ws.Range("B9").Value = 15
$ws.Range("D9").Value = 15
